# Apply the edits described by the commit:
# "Corrected tests to precreate containers into the container handler hash
#  to avoid attempts to access via solr."
#
# This updates the Container Profile Record No. (column M) for rows 6-8
# from 54555 -> 54556, and the Top Container Record No. (column J) for
# rows 9-10 from the old solr-looked-up ids (154691 / 154692) to the
# precreated container id (4). It also updates the active selection to
# match the author's final cursor position (M7:M8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testTopLinkerUpload")

# Column M ("Container Profile Record No.") for data rows 6, 7, 8
$ws.Cells.Item(6, 13).Value = 54556
$ws.Cells.Item(7, 13).Value = 54556
$ws.Cells.Item(8, 13).Value = 54556

# Column J ("Top Container Record No.") for the child rows 9, 10
$ws.Cells.Item(9, 10).Value = 4
$ws.Cells.Item(10, 10).Value = 4

# Move/restore the selection to match the saved cursor position
$ws.Range("M7:M8").Select()
